$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B ("Distribution Tpe" etc. shift right by one).
#    Excel copies column A's formatting (width/style) into the freshly inserted column,
#    matching the author's "insert column" edit.
$ws.Columns("B:B").Insert()

# 2. New header cell B1: "Dataset\nStandardised?" - same look as A1 (bold, centered, wrapped).
$ws.Range("B1").Value = "Dataset`nStandardised?"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108
$ws.Range("B1").WrapText = $true

# 3. New column B data (rows 2-8): "no" for all existing rows, centered like column A.
$bValues = @("no","no","no","no","no","no","no")
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $r = 2 + $i
    $cell = $ws.Range("B$r")
    $cell.Value = $bValues[$i]
    $cell.HorizontalAlignment = -4108
}

# 4. Column B width to match column A.
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# 5. Append the two new experiment rows (9 and 10).
function Set-Centered($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = $val
    $c.HorizontalAlignment = -4108
}

# Row 9: no / yes / kernel / no / 0.31 / 0.32 / 385.67 / (blank) / yes
Set-Centered "A9" "no"
Set-Centered "B9" "yes"
$ws.Range("C9").Value = "kernel"
Set-Centered "D9" "no"
$ws.Range("E9").Value = 0.31
$ws.Range("F9").Value = 0.32
$ws.Range("G9").Value = 385.67
Set-Centered "I9" "yes"

# Row 10: no / yes / kernel / yes / 0.25 / 0.26 / 4419.38 / 0.11377 / yes
Set-Centered "A10" "no"
Set-Centered "B10" "yes"
$ws.Range("C10").Value = "kernel"
Set-Centered "D10" "yes"
$ws.Range("E10").Value = 0.25
$ws.Range("F10").Value = 0.26
$ws.Range("G10").Value = 4419.38
$ws.Range("H10").Value = 0.11377
Set-Centered "I10" "yes"

# 6. Update the active selection to mirror the author's last-saved cursor position.
$ws.Range("I12").Select()
